# Generate Report for Handback
# Updates status text, timestamps, and clears stale "version not latest" error
# message for the zh-cn and de-de handback rows, plus widens/narrows a few
# report columns to better fit the new content.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# This shared string is used on Overview!E2, Overview!F2, zh-cn!C2, de-de!C2
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"

# Latest Handback DateTime (column K) refreshed to the new handback time
$wsZhCn.Range("K2").Value = "2016-08-29 16:54:30"
$wsDeDe.Range("K2").Value = "2016-08-29 16:54:37"

# Error Detail (column P) cleared now that the handback is in sync
$wsZhCn.Range("P2").Value = ""
$wsDeDe.Range("P2").Value = ""

# Column width adjustments to better accommodate the new report content
# (ColumnWidth is expressed in characters of the Normal style font; Excel
# stores width = ColumnWidth + 5/6 internally, snapped to whole pixels)
$wsOverview.Columns.Item(5).ColumnWidth = 29.1443713960194
$wsOverview.Columns.Item(6).ColumnWidth = 29.1443713960194

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1443713960194
$wsZhCn.Columns.Item(16).ColumnWidth = 12.9137195405506

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1443713960194
$wsDeDe.Columns.Item(16).ColumnWidth = 12.9137195405506
